# Applies the MT_10.xlsx edit described by the diff:
#  - Sheet2 (the active sheet): replace the plain numeric values in B2:F4
#    with "=1/<value>" formulas (reciprocals), keep the same cached result.
#  - Move the sheet's selection to F5.
#  - Touch the page setup (portrait orientation) so a <pageSetup/> element
#    is emitted, matching the upload's print-setup metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Row 2 (thread_id 23)
$ws.Range("B2").Formula = "=1/0.011"
$ws.Range("C2").Formula = "=1/0.299426117765301"
$ws.Range("D2").Formula = "=1/7.61700886175144"
$ws.Range("E2").Formula = "=1/0.014"
$ws.Range("F2").Formula = "=1/0.291590466236467"

# Row 3 (thread_id 24)
$ws.Range("B3").Formula = "=1/0.0790189850605536"
$ws.Range("C3").Formula = "=1/0.347690954728477"
$ws.Range("D3").Formula = "=1/0.0397994974842649"
$ws.Range("E3").Formula = "=1/0.050635955604689"
$ws.Range("F3").Formula = "=1/0.127279220613579"

# Row 4 (thread_id 25)
$ws.Range("B4").Formula = "=1/0.0807527089328899"
$ws.Range("C4").Formula = "=1/0.390037177715151"
$ws.Range("D4").Formula = "=1/7.61075817510975"
$ws.Range("E4").Formula = "=1/0.0451220566907141"
$ws.Range("F4").Formula = "=1/0.165786006647123"

# Move the active selection to F5, as in the saved workbook.
[void]$ws.Range("F5").Select()

# Touch page setup (portrait) so a pageSetup element is written for the sheet.
$ws.PageSetup.Orientation = 1
